# "Updated file to show our actual velocity in sprint 2."
#
# Sheet1 holds the Planned vs Actual velocity table:
#        B (Iteration 1)   C (Iteration 2)
# Planned      30                40
# Actual       37              <blank>   -> fill in with the real value (40)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Record the Sprint 2 (Iteration 2) actual velocity.
$ws.Range("C3").Value = 40

# Leave the selection where the author left it when saving.
$ws.Range("C5").Select()
